$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "saco acentos de los TCs" - update the test-case (TC) user and claim number
# in the sample row so they no longer carry the old values.
$ws.Range("D2").Value = "apellegrini"
$ws.Range("F2").Value = "1120194100405"

# Update the active selection left on the sheet.
$ws.Activate()
$ws.Range("D3").Select()
